$d = $word.ActiveDocument

# Target 1
$rng = $d.Content
$rng.Find.Text = 'Enable the students with the knowledge needed to use the concepts of psychology in the workplace.Promote the recognition of workers'' needs, both in its area of operation as in human relationships that have in the company.Promote strategies for the development of transversal skills necessary for the future engineer: effective communication, teamwork, commitment, initiative, responsibility, ethics, among others.'
$found = $rng.Find.Execute()
if (-not $found) { Write-Host "NOT FOUND: Target 1" }
$newText = 'Enable the students with the knowledge needed to use the concepts of psychology in the workplace.' + [char]11 + 'Promote the recognition of workers'' needs, both in its area of operation as in human relationships that have in the company.' + [char]11 + 'Promote strategies for the development of transversal skills necessary for the future engineer: effective communication, teamwork, commitment, initiative, responsibility, ethics, among others.'
$rng.Text = $newText

# Target 2
$rng = $d.Content
$rng.Find.Text = 'Possibilitar aos alunos da disciplina o conhecimento necessário para a utilização dos conceitos da psicologia em ambiente de trabalho.Favorecer o reconhecimento das necessidades dos trabalhadores tanto na sua área de atuação quanto nos relacionamentos humanos que terá na empresa.Promover estratégias para o desenvolvimento de competências transversais necessárias ao futuro engenheiro: comunicação eficaz, trabalho em equipe, comprometimento, iniciativa, responsabilidade, ética, entre outras.'
$found = $rng.Find.Execute()
if (-not $found) { Write-Host "NOT FOUND: Target 2" }
$newText = 'Possibilitar aos alunos da disciplina o conhecimento necessário para a utilização dos conceitos da psicologia em ambiente de trabalho.' + [char]11 + 'Favorecer o reconhecimento das necessidades dos trabalhadores tanto na sua área de atuação quanto nos relacionamentos humanos que terá na empresa.' + [char]11 + 'Promover estratégias para o desenvolvimento de competências transversais necessárias ao futuro engenheiro: comunicação eficaz, trabalho em equipe, comprometimento, iniciativa, responsabilidade, ética, entre outras.'
$rng.Text = $newText

# Target 3
$rng = $d.Content
$rng.Find.Text = '1.Introdução: conceituar psicologia como ciência e como aplicação; a psicologia aplicada ao trabalho. A psicologia nas relações humanas no trabalho.2.Conceito de Comunicação: sistemas, funções, axiomas da comunicação humana. Processos de comunicação e o convívio sócio-comunicacional na empresa.3.Relações Humanas no Trabalho: relações humanas em grupos; como participar de um grupo de trabalho, trabalho em equipe, dinâmicas grupais.4.Psicologia nas Organizações de Trabalho: conceitos de organização e de trabalho. Organização e trabalho e sua importância na saúde mental e produtividade do trabalhador: estresse, síndrome de burnout, síndrome de Karoshi; L.E.R.; qualidade de vida; assédios sexual e moral no ambiente de trabalho.5.Recrutamento e Seleção: recrutamento e seleção de pessoal; colocação e acompanhamento; avaliação de desempenho; medidas de avaliação e sua importância na seleção; experiências práticas em sala de aula como facilitadoras do processo de seleção.6.Motivação: as necessidades básicas e psicológicas do ser humano; motivação e conflitos; fatores esquecidos como motivadores na empresa: inveja, ciúme, medo, abuso de poder. Avaliação de motivação.7. Liderança: definição, teorias e desenvolvimento de lideranças8. Treinamento e Desenvolvimento: definição, diferenciação, etapas, dificuldades9. Avaliação de desempenho: definição, tipos, periodicidade, importância'
$found = $rng.Find.Execute()
if (-not $found) { Write-Host "NOT FOUND: Target 3" }
$newText = '1.Introdução: conceituar psicologia como ciência e como aplicação; a psicologia aplicada ao trabalho. A psicologia nas relações humanas no trabalho.' + [char]11 + '2.Conceito de Comunicação: sistemas, funções, axiomas da comunicação humana. Processos de comunicação e o convívio sócio-comunicacional na empresa.' + [char]11 + '3.Relações Humanas no Trabalho: relações humanas em grupos; como participar de um grupo de trabalho, trabalho em equipe, dinâmicas grupais.' + [char]11 + '4.Psicologia nas Organizações de Trabalho: conceitos de organização e de trabalho. Organização e trabalho e sua importância na saúde mental e produtividade do trabalhador: estresse, síndrome de burnout, síndrome de Karoshi; L.E.R.; qualidade de vida; assédios sexual e moral no ambiente de trabalho.' + [char]11 + '5.Recrutamento e Seleção: recrutamento e seleção de pessoal; colocação e acompanhamento; avaliação de desempenho; medidas de avaliação e sua importância na seleção; experiências práticas em sala de aula como facilitadoras do processo de seleção.' + [char]11 + '6.Motivação: as necessidades básicas e psicológicas do ser humano; motivação e conflitos; fatores esquecidos como motivadores na empresa: inveja, ciúme, medo, abuso de poder. Avaliação de motivação.' + [char]11 + '7. Liderança: definição, teorias e desenvolvimento de lideranças' + [char]11 + '8. Treinamento e Desenvolvimento: definição, diferenciação, etapas, dificuldades' + [char]11 + '9. Avaliação de desempenho: definição, tipos, periodicidade, importância'
$rng.Text = $newText

# Target 4
$rng = $d.Content
$rng.Find.Text = 'Serão aplicadas provas dissertativas com estudo de caso e situações, para levar os alunos à maior reflexão sobre a utilização dos conceitos aprendidos para o futuro engenheiro em seu trabalho cotidiano em empresas.Será solicitada a realização de atividades variadas (avaliação processual) sobre cada tema, com ênfase no desenvolvimento das habilidades transversais. Tais atividades poderão ser: apresentações, elaboração de folder, mapas conceituais e pitch para processo seletivo, relatório e leitura ativa, entrevista com trabalhadores, etc.Será realizada uma atividade extensionista de confecção de currículo para a comunidade. Para tanto, os alunos deverão participar de algum evento que ocorra durante o semestre (feiras de ciências, feira de profissões, reunião com calouros de outras instituições). Caberá ao aluno o planejamento e a execução da atividade, com supervisão da profa.'
$found = $rng.Find.Execute()
if (-not $found) { Write-Host "NOT FOUND: Target 4" }
$newText = 'Serão aplicadas provas dissertativas com estudo de caso e situações, para levar os alunos à maior reflexão sobre a utilização dos conceitos aprendidos para o futuro engenheiro em seu trabalho cotidiano em empresas.' + [char]11 + 'Será solicitada a realização de atividades variadas (avaliação processual) sobre cada tema, com ênfase no desenvolvimento das habilidades transversais. Tais atividades poderão ser: apresentações, elaboração de folder, mapas conceituais e pitch para processo seletivo, relatório e leitura ativa, entrevista com trabalhadores, etc.' + [char]11 + 'Será realizada uma atividade extensionista de confecção de currículo para a comunidade. Para tanto, os alunos deverão participar de algum evento que ocorra durante o semestre (feiras de ciências, feira de profissões, reunião com calouros de outras instituições). Caberá ao aluno o planejamento e a execução da atividade, com supervisão da profa.'
$rng.Text = $newText

# Target 5
$rng = $d.Content
$rng.Find.Text = 'P1 = soma de todas as atividades até 8,0 + relatório do desenvolvimento da atividade de extensão.P2 = prova dissertativa aplicada no final do semestre.Média = (P1+ P2)/ 2.'
$found = $rng.Find.Execute()
if (-not $found) { Write-Host "NOT FOUND: Target 5" }
$newText = 'P1 = soma de todas as atividades até 8,0 + relatório do desenvolvimento da atividade de extensão.' + [char]11 + 'P2 = prova dissertativa aplicada no final do semestre.' + [char]11 + 'Média = (P1+ P2)/ 2.'
$rng.Text = $newText

# Target 6
$rng = $d.Content
$rng.Find.Text = '1.ROBBINS, S. P. Fundamentos do Comportamento Organizacional. 8ª. Ed., São Paulo: Pearson Prentice Hall, 2009.2.CHIAVENATO, I. Gestão de Pessoas: o novo papel da Gestão do Talento Humano. 5ª.ed., São Paulo: Atlas, 2020.3.SHILD, M (trad) A arte de dar feedback. Harvard Business Review. Rio de Janeiro: Sextante, 2019.4.BROWNIE, S. HR on Purpose: Developing Deliberate People Passion — Steve Browne. Alexandria: Society for Human Resource Management, 2017.5.BERGAMINI, C.W. Motivação nas Organizações.7ª.Ed. São Paulo: Atlas, 2018.6.ROSENBERG, M.B. Comunicação Não-Violenta. Rio de Janeiro: Sextante, 2018.'
$found = $rng.Find.Execute()
if (-not $found) { Write-Host "NOT FOUND: Target 6" }
$newText = '1.ROBBINS, S. P. Fundamentos do Comportamento Organizacional. 8ª. Ed., São Paulo: Pearson Prentice Hall, 2009.' + [char]11 + '2.CHIAVENATO, I. Gestão de Pessoas: o novo papel da Gestão do Talento Humano. 5ª.ed., São Paulo: Atlas, 2020.' + [char]11 + '3.SHILD, M (trad) A arte de dar feedback. Harvard Business Review. Rio de Janeiro: Sextante, 2019.' + [char]11 + '4.BROWNIE, S. HR on Purpose: Developing Deliberate People Passion — Steve Browne. Alexandria: Society for Human Resource Management, 2017.' + [char]11 + '5.BERGAMINI, C.W. Motivação nas Organizações.7ª.Ed. São Paulo: Atlas, 2018.' + [char]11 + '6.ROSENBERG, M.B. Comunicação Não-Violenta. Rio de Janeiro: Sextante, 2018.'
$rng.Text = $newText

# Target 7
$rng = $d.Content
$rng.Find.Text = 'Introduction : conceptualize psychology as science and application; psychology applied to work. The psychology of human relations at work. Concept of Communication : Systems, functions , axioms of human communication. Communication processes and the social and communicative interaction in the company. Human relations at work: the role of masks in human interaction; human relations in groups; how to be a part of a workgroup, teamwork and group dynamics Psychology in Work Organizations : Organization concepts and work. Organization and work and their importance in mental health and worker productivity: stress, burnout , Karoshi syndrome ; L.E.R .; quality of life; sexual and moral harassment in the workplace ; alcohol and drugs at work; mental disorders in the company. Recruitment and Selection: recruitment and selection of personnel ; placement and monitoring; performance evaluation; training and education; evaluation measures and their importance in the selection ; practical experiences in the classroom as facilitators in the selection process . Motivation : the basic and psychological needs of human beings; motivation and conflicts; forgotten factors as motivators in the company : envy, jealousy , fear, abuse of power . Motivation Assessment.- Leadership: definition, theories and leadership development- Training and Development: definition, differentiation, stages, difficulties- Performance evaluation: definition, types, frequency, importance'
$found = $rng.Find.Execute()
if (-not $found) { Write-Host "NOT FOUND: Target 7" }
$newText = 'Introduction : conceptualize psychology as science and application; psychology applied to work. The psychology of human relations at work.' + [char]11 + ' Concept of Communication : Systems, functions , axioms of human communication. Communication processes and the social and communicative interaction in the company.' + [char]11 + ' Human relations at work: the role of masks in human interaction; human relations in groups; how to be a part of a workgroup, teamwork and group dynamics' + [char]11 + ' Psychology in Work Organizations : Organization concepts and work. Organization and work and their importance in mental health and worker productivity: stress, burnout , Karoshi syndrome ; L.E.R .; quality of life; sexual and moral harassment in the workplace ; alcohol and drugs at work; mental disorders in the company.' + [char]11 + ' Recruitment and Selection: recruitment and selection of personnel ; placement and monitoring; performance evaluation; training and education; evaluation measures and their importance in the selection ; practical experiences in the classroom as facilitators in the selection process .' + [char]11 + ' Motivation : the basic and psychological needs of human beings; motivation and conflicts; forgotten factors as motivators in the company : envy, jealousy , fear, abuse of power . Motivation Assessment.' + [char]11 + '- Leadership: definition, theories and leadership development' + [char]11 + '- Training and Development: definition, differentiation, stages, difficulties' + [char]11 + '- Performance evaluation: definition, types, frequency, importance'
$rng.Text = $newText

# Target 8
$rng = $d.Content
$rng.Find.Text = 'Espera-se que o aluno de graduação consiga aplicar conhecimento aprendido na disciplina LOB 1031, ajudando o publico-alvo a aprimorar e ampliar sua possibilidade de empregabilidade, uma vez que o currículum vitae é exigência em praticamente todos os processos de recrutamento e seleção.Espera-se que o grupo social atendido tenha mais condições de aumentar suas possibilidades de participação em processos seletivos.'
$found = $rng.Find.Execute()
if (-not $found) { Write-Host "NOT FOUND: Target 8" }
$newText = 'Espera-se que o aluno de graduação consiga aplicar conhecimento aprendido na disciplina LOB 1031, ajudando o publico-alvo a aprimorar e ampliar sua possibilidade de empregabilidade, uma vez que o currículum vitae é exigência em praticamente todos os processos de recrutamento e seleção.' + [char]11 + 'Espera-se que o grupo social atendido tenha mais condições de aumentar suas possibilidades de participação em processos seletivos.'
$rng.Text = $newText

# Target 9
$rng = $d.Content
$rng.Find.Text = 'O aluno deverá, no último mês de aula do semestre:- fazer uma pesquisa aprofundada sobre formas e formatos de curriculum vitae- fazer contato com entidades ou organizadores de evento para oferecer e acertar detalhes da aplicação da atividade- reunir-se com o participante e desenvolver a atividade proposta: explicar do que se trata um CV, coletar as informações necessárias, elaborar conjuntamente o CV, fornecer orientações básicas para apresentação em entrevista de emprego, aplicar o questionário de satisfação'
$found = $rng.Find.Execute()
if (-not $found) { Write-Host "NOT FOUND: Target 9" }
$newText = 'O aluno deverá, no último mês de aula do semestre:' + [char]11 + '- fazer uma pesquisa aprofundada sobre formas e formatos de curriculum vitae' + [char]11 + '- fazer contato com entidades ou organizadores de evento para oferecer e acertar detalhes da aplicação da atividade' + [char]11 + '- reunir-se com o participante e desenvolver a atividade proposta: explicar do que se trata um CV, coletar as informações necessárias, elaborar conjuntamente o CV, fornecer orientações básicas para apresentação em entrevista de emprego, aplicar o questionário de satisfação'
$rng.Text = $newText
